$wb = $excel.ActiveWorkbook

# ---- Sheet: YDS (long per-play running logs, R and P rows) ----
$wsYDS = $wb.Worksheets.Item("YDS")
$wsYDS.Range("B2").Value2 = '3 6 0 0 3 23 6 1 7 6 -2 4 14 2 1 -1 5 4 1 3 5 -4 5 7 4 2 7 7 5 0 5 11 3 2 1 6 -1 5 3 -4 -7 9 5 5 1 11 2 5 3 2 1 2 4 0 1 0 9 0 -1 6 9 0 5 2 4 2 0 3 7 0 5 13 1 2 4 11 34 8 23 8 1 4 -1 2 5 11 10 7 8 21 5 -1 5 5 1 6 3 1 -3 2 2 2 2 6 3 1 2 -1 4 3 7 2 0 7 2 0 2 2 -1 6 9 2 2 2 1 -1 3 6 3 6 7 2 25 -6 0 3 4 2 3 4 9 1 -4 0 2 0 5 0 1 -1 7 0 1 1 8 1 3 5 3 5 7 7 4 1 -3 12 12 6 1 2 6 9 2 6 2 1 1 12 -3 4 0 12 4 7 2 4 5 2 1 16 1 -1 6 3 1 2 0 3 1 3 5 1 3 1 1 14 8 2 11 4 1 6 0 12 1 15 -1 2 8 7 4 1 4 3 5 10 0 1 7 7 0 3 2 8 0 2 19 2 -5 0 1 2 -4 4 2 15 2 1 8 4 -2 1 9 1 -1 0 3 2 5 1 3 3 0 10 0 -1 -1 4 1 3 6 6 0 3 -2 2 9 2 4 0 9 0 3 5 9 4 3 4 11 3 0 1 5 4 0 5 8 0 -5 7 1 1 4 6 5 0 1 6 1 4 4 11 1 2 3 -1 3 -1 3 5 7 9 1 4 8 8 23 4 5 1 3 4 -1 5 4 0 1 3 2 1 14 0 2 2 6 0 3 5 2 1 11 46 2 7 5 1 1 2 9 3 11 7 5 3 7 2 2 4 2 1 3 4 13 1 11 2 -2 55 2 -2 4 6 15 9 2 4 17 5 -2 3 4 1 2 -1 13 1 3 6 4 7 2 2 5 10 19 1 5 5 9 12 -1 7 1 2 8 2 4 0 3 0 5 2 3 2 3 1 8 -1 -3 4 9 5 1 5 10 3 3 5 2 4 3 2 6 4 -2 6 27 1 -1 3 11 11 4 3 2 8 1 1 5 2 7 -1 2 2 2 7 1 0 4 3 9 4 11 5 7 0 1 4 0 1 8 1 2 4 0 7 10 -1 9 4 0 6 5 7 11 1 -1 4 9 8 2 3 -2 10 4 4 1 1 0 9 8 5 2 3 0 2 4 6 6 14 4 3 5 1 7 6 3 4 6 2 4 5 18 4 2 1 -2 6 4 2 9 0 -1 -1 9 0 0 5 1 0 0 10 7 21 3 46 3 -3 -5 2 2 1 -1 -1 -1 1 6 8 0 -4 6 7 4 9 4 5 -1 2 1 2 0 3 3 7 4 3 0 21 2 2 4 0 1 2 5 0 2 9 2 6 4 14 4 5 1 3 -2 -3 1 -2 10 -1 11 2 5 4 8 -2 18 5 2 6 4 -5 2 1 19 4 0 1 3 7 20 1 0 1 25 3 7 2 3 6 3 8 3 1 6 4 5 9 8 9 0 9 5 1 8 9 4 -1 0 1 32 1 1 2 3 13 2 3 -2 -1 4 0 2 2 7 1 -5 6 8 6 2 6 17 2 5 3 4 0 10 3 3 8 -3 5 10 -2 8 3 3 3 2 4 6 4 0 2 7 1 6 1 2 4 6 2 0 17 1 0 4 3 6 0 2 4 5 5 6 5 1 -1 5 8 4 3 5 22 2 12 1 8 2 4 1 1 3 8 -2 19 3 1 1 9 1 1 1 3 12 2 2 2 8 2 2 2 13 3 6 6 3 1 4 3 3 1 3 4 5 4 1 0 0 0'
$wsYDS.Range("B3").Value2 = '2 11 7 4 2 13 6 1 5 9 4 11 14 15 19 4 8 15 10 12 5 6 8 13 4 5 7 18 23 5 11 6 18 6 14 5 15 12 5 5 9 4 5 6 3 6 5 4 10 9 8 2 5 3 9 9 15 19 16 5 1 2 7 5 7 13 17 25 4 12 10 1 14 18 42 11 1 8 12 0 11 3 15 20 0 9 15 7 16 30 23 14 8 9 14 12 13 3 9 24 0 5 16 15 -3 23 14 26 12 -2 14 6 7 10 4 16 8 11 3 5 12 9 9 7 22 12 12 17 23 67 4 6 5 5 19 13 9 7 7 8 10 20 11 7 11 8 2 2 6 17 -2 29 14 13 10 7 22 25 4 4 3 7 27 5 23 13 9 11 16 6 10 16 5 6 9 14 16 10 16 12 15 25 9 9 0 10 24 2 6 6 7 8 13 2 13 3 7 12 8 12 13 9 13 7 4 18 6 15 8 -4 5 20 9 54 5 7 13 8 13 5 9 16 3 13 0 5 3 18 9 5 6 6 5 13 -5 11 12 12 19 8 9 25 5 15 11 2 4 5 8 5 11 14 9 16 11 1 5 1 10 2 4 13 7 14 4 5 18 16 1 8 3 12 7 4 5 72 15 11 7 9 8 -3 5 6 2 10 11 16 6 7 21 11 5 7 9 10 0 6 10 5 4 9 3 1 7 2 5 4 16 9 5 5 9 30 14 8 14 4 10 28 8 15 6 8 17 8 5 6 9 7 5 33 30 -6 9 6 12 7 14 42 20 14 22 7 11 20 4 1 8 21 8 7 3 5 3 7 12 17 13 2 14 2 9 50 4 28 5 14 15 18 5 6 2 7 12 10 3 12 1 2 14 8 7 14 14 22 5 10 15 42 7 12 15 7 17 7 19 34 14 8 14 5 9 4 7 33 -2 15 18 0 10 21 0 44 22 25 12 1 11 8 4 8 12 7 31 12 15 9 7 25 2 7 9 6 7 2 10 4 7 5 70 28 14 5 8 3 10 2 7 19 9 8 8 5 21 14 24 5 5 1 9 4 3 7 9 34 8 7 23 11 53 14 2 7 11 23 6 5 9 55 7 4 6 13 26 7 17 21 32 7 82 25 2 11 11 27 6 3 13 0 3 9 16 4 26 2 46 9 21 8 54 8 10 11 6 21 10 19 11 11 9 20 15 19 16 9 9 5 11 16 6 12 5 5 12 8 11 9 3 3 4 11 11 11 9 11 6 3 6 3 2 7 14 9 14 3 7 6 4 17 7 3 7 9 9 12 6 5 9 5 16 4 12 32 11 16 7 15 16 6 0 -4 0 -7 16 12 19 5 12 21 3 9 7 32 16 29 5 7 0 5 9 20 12 6 9 11 12 16 14 25 14 10 14 20 0 10 11 16 9 12 -4 11 12 19 3 5 27 4 13 5 12 22 29 32 26 23 6 2 6 17 2 2 9 8 4 19 3 56 4 4 15 10 18 20 3 11 16 5 17 3 68 20 17 9 11 15 3 52 1 19 3 18 9 3 2 5 14 13 0 12 4 22 3 10 20 4 52 11 4 3 10 72 8 8 9 17 18 6 7 6 4 10 5 69 22 6 4 39 10 5 2 5 9 2 35 30'
$wsYDS.Range("C2").Value2 = '4 3 2 3 2 2 13 4 1 13 3 0 7 2 6 7 5 3 6 26 1 2 3 4 2 4 8 4 5 2 8 4 -1 -1 0 2 7 5 -2 4 5 1 3 6 8 6 11 1 13 -1 6 6 2 1 3 6 11 2 22 0 -1 1 4 4 11 10 6 4 26 33 1 10 1 11 14 0 1 4 9 12 3 2 3 3 4 24 1 1 19 5 2 2 8 4 -6 3 -3 4 7 0 6 9 7 7 2 2 2 1 2 4 8 4 0 8 2 2 0 14 3 3 5 2 7 40 4 13 4 3 -2 42 3 4 1 1 3 -4 4 0 6 5 3 34 0 3 10 3 14 3 20 4 3 3 0 1 1 2 -1 1 7 7 3 21 17 -2 2 1 0 6 11 5 12 -4 15 2 3 -1 5 -3 3 0 5 2 3 0 3 11 2 3 0 15 4 0 6 5 6 2 9 13 8 11 3 6 4 3 3 16 4 8 21 15 2 9 5 7 8 17 8 5 1 7 1 16 0 -6 -3 -5 1 1 0 11 1 3 0 -1 2 2 8 14 6 2 6 3 1 3 6 5 4 1 19 16 3 2 1 2 6 2 4 16 10 2 3 13 8 -1 -5 1 3 11 5 4 2 1 0 4 5 3 13 2 1 4 -1 2 4 5 18 4 7 2 4 2 5 2 2 -1 4 7 5 -2 3 -1 3 3 4 2 8 3 3 4 2 4 11 4 0 1 1 3 5 0 4 0 0 14 1 4 2 6 1 26 7 1 4 2 2 5 3 -1 5 0 4 1 15 0 1 1 17 2 -1 4 7 26 -2 2 4 3 2 2 0 5 3 6 0 13 2 -6 3 4 7 2 29 13 3 -2 6 -2 -3 1 3 2 2 4 10 6 2 5 5 38 4 48 14 2 10 0 14 4 3 1 13 7 2 22 2 12 3 5 6 3 6 5 5 5 8 22 9 0 3 2 5 3 2 3 18 6 27 20 4 4 7 4 4 11 4 8 72 16 6 1 16 2 2 5 2 0 1 2 0 1 1 3 1 2 17 2 7 -3 4 1 13 0 8 -1 6 5 2 1 0 -1 0 0 2 2 3 7 8 6 11 14 5 0 -1 0 6 6 5 11 7 2 0 -3 6 1 1 4 4 -1 1 10 6 4 4 2 4 -3 5 1 1 3 20 0 2 2 1 2 1 11 3 6 2 4 11 2 11 12 6 1 3 5 5 20 1 3 0 7 2 3 0 8 0 2 4 1 4 2 14 11 3 0 2 1 3 8 1 5 7 17 -1 5 4 1 1 57 0 -1 -4 -1 3 1 0 3 8 1 6 2 1 1 -5 1 0 9 5 -2 1 2 10 -5 4 -2 -1 12 3 10 3 7 1 16 2 1 7 7 2 10 7 2 5 6 0 6 8 6 2 5 9 -1 1 8 4 4 3 2 5 5 18 2 -1 1 4 2 7 1 11 22 0 1 2 3 0 1 70 7 1 -2 13 10 2 3 2 -1 -1 -4 6 4 3 3 6 4 5 1 7 4 1 4 3 2 4 5 0 18 -4 7 2 3 1 2 2 4 5 4 3 5 0 4 3 -1 13 3 2 0 -1 6 7 3 9 5 0 2 5 4 1 15 1 5 6 1 3 -4 7 2 -1 0 1 4 4 9 12 2 27 1 7 2 3 -3 2 5 3 3 1 8 6 2 2 4 2 4 4 3 5 -1 5 2 2 1 3 2 6 7 2 2 5 6 14 4 5 7 1 0 -1 7 7 13 1 -2 10 1 2 8 2 6 1 0 7 -2 3 -5 3 4 2 2 -7 1 1 2 23 -1 6 2 1 8 10 6 8 1 14 3 17 1 19 8 0 5 8 11'
$wsYDS.Range("C3").Value2 = '10 4 13 15 7 7 22 10 33 7 3 16 6 8 1 10 10 3 3 8 9 26 11 9 43 7 16 18 9 13 6 0 11 6 14 7 9 7 12 2 6 3 0 12 16 3 10 8 29 9 10 0 3 11 3 -2 11 4 5 0 4 30 6 22 7 4 8 7 14 11 14 20 9 22 7 28 9 7 13 13 5 4 7 33 2 18 26 9 10 7 16 9 9 14 5 14 2 7 2 9 12 9 12 3 7 5 8 30 22 -2 7 5 17 15 55 22 13 21 10 7 6 20 11 1 20 17 5 5 8 10 13 10 14 2 11 20 6 28 -6 10 11 3 35 8 4 12 6 4 6 21 21 16 20 19 12 8 13 8 30 24 10 5 27 24 9 23 19 0 6 9 11 9 10 26 12 10 12 13 1 10 -6 4 46 16 12 14 4 15 24 8 18 24 17 3 6 11 8 14 12 11 6 9 8 5 10 0 8 42 2 3 4 11 7 14 3 25 12 11 7 7 6 6 4 53 6 12 9 11 7 10 18 16 5 7 5 44 2 -9 13 -1 8 15 5 4 8 6 15 8 16 20 7 9 2 11 5 5 17 14 35 5 22 11 8 20 4 8 9 9 18 11 5 4 2 20 13 14 12 11 11 14 7 15 7 32 4 -2 3 0 3 9 12 37 23 8 8 4 0 9 21 4 14 6 8 9 3 7 50 2 12 17 7 26 25 11 6 7 14 13 -3 4 2 29 13 0 22 15 12 5 10 43 5 9 18 5 2 9 8 9 6 12 11 9 -1 2 9 4 17 11 25 5 5 6 20 0 9 10 10 5 6 9 8 34 10 7 27 0 6 14 7 3 5 12 10 17 -5 3 3 11 2 36 10 5 11 10 13 21 5 0 5 3 0 6 3 7 3 16 12 6 6 3 18 4 6 15 21 16 4 24 0 6 8 3 5 8 7 5 8 11 17 5 29 9 3 8 8 8 17 7 12 10 7 7 52 7 4 27 5 20 -2 7 11 5 15 17 -4 34 14 0 12 13 24 -1 3 15 24 5 22 5 7 0 8 8 2 2 59 20 20 15 0 7 10 4 6 5 3 5 -2 5 4 33 2 9 5 3 15 7 4 3 7 10 3 7 8 10 14 15 24 6 23 19 20 18 7 1 35 39 25 3 13 13 -1 10 12 7 11 9 30 8 7 6 3 15 3 6 7 12 6 20 6 5 2 26 11 12 15 9 4 8 13 14 15 9 8 9 10 24 14 23 12 6 8 8 19 13 3 23 3 11 60 15 6 5 22 7 -2 21 26 18 -2 25 6 8 6 21 4 1 9 2 25 31 19 7 8 20 9 8 8 11 41 0 11 7 11 2 15 3 10 12 2 30 9 16 13 9 6 8 15 6 5 10 11 -1 10 -1 41 5 10 47 6 7 44 23 7 4 7 12 7 12 9 7 6 -3 1 18 33 2 7 5 4 2 10 21 6 1 4 3 22 20 16 10 13 8 9 5 17 19 8 25 21 3 9 12 5 8 7 28 -1 10 3 11 -1 8 9 11 8 25 2 24 24 5 11 17 5 15 6 4 9 -1 16 12 10 -1 18 28 11 11 4 7 12 17 25 18 14 4 6 7 14 5 8 6 1 29 7 17 24 3 3 53 7 3 2 4 27 7 8 5 11 6 5 6 4 8 1 4'

# ---- Sheet: ST (special teams running logs, D and RA rows) ----
$wsST = $wb.Worksheets.Item("ST")
$wsST.Range("B4").Value2 = '66 64 58 67 63 64 45 64 66 66 51 63 58 60 62 68 60 66 59 70 55 63 69 56 62 67 62 62 66 64 67 64 39 62 68 62 66 59 67 68 68 65 55 63 67 66 56 64 63 64 64 67 45'
$wsST.Range("B5").Value2 = '46 0 30 27 19 24 11 19 16 23 10 25 17 8 17 60 30 23 12 21 10 22 27 16 20 23 28 18 21 24 35 24 11 12 25 15 24 24 25 48 24 27 16 24 23 22 15 24 0 25 26 30 12'
$wsST.Range("D3").Value2 = '70 53 57 51 55 44 50 48 38 39 45 33 46 48 49 40 57 37 46 59 45 47 60 53 50 28 41 42 47 43 46 41 50 37 42 45 48 49 50 45 41 59 48 34 55 47 50 57 51 52 53 41 48 51 39 49 49 41 51 46 48 35 40 53 56 37 30 47 72 44 59 47 37 43 42 43 61 36 39 48 37 54 49 53 46 56 40 46 41 38 42 25 58 51 54 39 60 58 59 47 45 38 55 37 44 41 47 58 41 45 53 50 46 46 38 46 36 41 50 57 57 41 61 45 56 51 40 51 27'
$wsST.Range("D4").Value2 = '0 7 17 0 0 0 13 0 0 0 0 0 0 1 1 0 0 0 0 10 1 -1 0 10 7 0 0 3 0 0 15 0 0 0 8 -5 0 42 10 0 3 0 0 0 8 0 0 0 2 14 6 0 0 1 -1 0 0 0 11 0 8 0 0 14 0 0 0 6 0 0 0 0 0 0 0 0 12 0 0 14 0 0 0 0 0 0 0 1 0 0 8 0 11 13 14 0 16 21 0 0 0 0 16 0 -1 5 12 11 0 3 0 15 3 0 0 0 4 0 0 0 0 0 0 0 12 18 0 29 3'
$wsST.Range("D5").Value2 = '0 0 12 8 0 1 0 12 0 0 5 0 5 6 10 0 0 3 0'

# ---- Numeric season-total updates (Week 17 stats logged) ----
$wsOFF = $wb.Worksheets.Item("OFF")
$wsOFF.Range("B2").Value2 = 4
$wsOFF.Range("C2").Value2 = 387
$wsOFF.Range("D2").Value2 = 23
$wsOFF.Range("F2").Value2 = 106
$wsOFF.Range("G2").Value2 = 114
$wsOFF.Range("J2").Value2 = 55
$wsOFF.Range("L2").Value2 = 579
$wsOFF.Range("M2").Value2 = 405
$wsOFF.Range("O2").Value2 = 46
$wsOFF.Range("P2").Value2 = 30
$wsOFF.Range("Q2").Value2 = 1089
$wsOFF.Range("B3").Value2 = 18
$wsOFF.Range("C3").Value2 = 351
$wsOFF.Range("E3").Value2 = 73
$wsOFF.Range("F3").Value2 = 209
$wsOFF.Range("G3").Value2 = 73
$wsOFF.Range("H3").Value2 = 68
$wsOFF.Range("I3").Value2 = 106
$wsOFF.Range("J3").Value2 = 114
$wsOFF.Range("N3").Value2 = 49

$wsDEF = $wb.Worksheets.Item("DEF")
$wsDEF.Range("C2").Value2 = 381
$wsDEF.Range("D2").Value2 = 19
$wsDEF.Range("F2").Value2 = 135
$wsDEF.Range("G2").Value2 = 94
$wsDEF.Range("J2").Value2 = 72
$wsDEF.Range("L2").Value2 = 569
$wsDEF.Range("M2").Value2 = 378
$wsDEF.Range("O2").Value2 = 34
$wsDEF.Range("P2").Value2 = 18
$wsDEF.Range("Q2").Value2 = 1055
$wsDEF.Range("B3").Value2 = 19
$wsDEF.Range("C3").Value2 = 384
$wsDEF.Range("F3").Value2 = 221
$wsDEF.Range("G3").Value2 = 72
$wsDEF.Range("H3").Value2 = 57
$wsDEF.Range("I3").Value2 = 125
$wsDEF.Range("J3").Value2 = 109

$wsST = $wb.Worksheets.Item("ST")
$wsST.Range("B2").Value2 = 163
$wsST.Range("D2").Value2 = 129
$wsST.Range("F2").Value2 = 48
$wsST.Range("G2").Value2 = 46
$wsST.Range("J2").Value2 = 52
$wsST.Range("K2").Value2 = 50
$wsST.Range("L2").Value2 = 25
$wsST.Range("M2").Value2 = 18
$wsST.Range("B3").Value2 = 110

$wsTURNS = $wb.Worksheets.Item("TURNS")
$wsTURNS.Range("E3").Value2 = 14

$wsPEN = $wb.Worksheets.Item("PEN")
$wsPEN.Range("B2").Value2 = 35
$wsPEN.Range("B3").Value2 = 34
$wsPEN.Range("D4").Value2 = 17
